$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-33: update Price (D) and/or Volume (E) values in place ---
$ws.Cells.Item(2, 4).Value = "26.510.69"
$ws.Cells.Item(2, 5).Value = "  -0.05%  "
$ws.Cells.Item(3, 4).Value = "1.732.70"
$ws.Cells.Item(3, 5).Value = "  +0.29%  "
$ws.Cells.Item(4, 5).Value = "  +0.21%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "244.49"
$ws.Cells.Item(5, 5).Value = "  -0.33%  "
$ws.Cells.Item(6, 5).Value = "  +0.15%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4911"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.2630"
$ws.Cells.Item(8, 5).Value = "  -1.49%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.06185"
$ws.Cells.Item(10, 4).Value = "1.735.03"
$ws.Cells.Item(10, 5).Value = "  +0.51%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07017"
$ws.Cells.Item(11, 5).Value = "  -1.80%  "
$ws.Cells.Item(12, 5).Value = "  -0.73%  "
$ws.Cells.Item(13, 5).Value = "  +0.82%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.6018"
$ws.Cells.Item(14, 5).Value = "  -2.44%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "77.47"
$ws.Cells.Item(15, 5).Value = "  +0.42%  "
$ws.Cells.Item(16, 5).Value = "  +0.16%  "
$ws.Cells.Item(17, 4).Value = "26.515.45"
$ws.Cells.Item(17, 5).Value = "  +0.00%  "
$ws.Cells.Item(18, 5).Value = "  +0.27%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.000007125"
$ws.Cells.Item(19, 5).Value = "  +2.76%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "11.44"
$ws.Cells.Item(20, 5).Value = "  -1.99%  "
$ws.Cells.Item(21, 4).Value = "1.961.34"
$ws.Cells.Item(21, 5).Value = "  +0.72%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "4.498"
$ws.Cells.Item(22, 5).Value = "  -0.66%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "8.620"
$ws.Cells.Item(23, 5).Value = "  -3.76%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "5.190"
$ws.Cells.Item(24, 5).Value = "  -1.85%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "138.97"
$ws.Cells.Item(25, 5).Value = "  +1.80%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "15.25"
$ws.Cells.Item(26, 5).Value = "  -0.56%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "1.420"
$ws.Cells.Item(27, 5).Value = "  +1.21%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "106.68"
$ws.Cells.Item(28, 5).Value = "  -0.13%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.724"
$ws.Cells.Item(29, 5).Value = "  -3.86%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "3.986"
$ws.Cells.Item(30, 5).Value = "  +0.12%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.07974"
$ws.Cells.Item(31, 5).Value = "  -0.68%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.689"
$ws.Cells.Item(32, 5).Value = "  -0.44%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.04541"
$ws.Cells.Item(33, 5).Value = "  -0.53%  "

# --- Rows 34-51: Frax row removed, list shifts up by one, Aave appended at the end ---
$ws.Cells.Item(34, 2).Value = "HuobiToken"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "2.617"
$ws.Cells.Item(34, 5).Value = "  +0.09%  "
$ws.Cells.Item(35, 2).Value = "ARBITRUM"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.007"
$ws.Cells.Item(35, 5).Value = "  +1.60%  "
$ws.Cells.Item(36, 2).Value = "ImmutableX"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.6257"
$ws.Cells.Item(36, 5).Value = "  -1.84%  "
$ws.Cells.Item(37, 2).Value = "TrustWalletToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.9057"
$ws.Cells.Item(37, 5).Value = "  -2.44%  "
$ws.Cells.Item(38, 2).Value = "RenderToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "2.004"
$ws.Cells.Item(38, 5).Value = "  -4.21%  "
$ws.Cells.Item(39, 2).Value = "MXToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "2.409"
$ws.Cells.Item(39, 5).Value = "  -0.19%  "
$ws.Cells.Item(40, 2).Value = "PaxDollar"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "1.002"
$ws.Cells.Item(40, 5).Value = "  -0.34%  "
$ws.Cells.Item(41, 2).Value = "VeChain"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.01489"
$ws.Cells.Item(41, 5).Value = "  -0.88%  "
$ws.Cells.Item(42, 2).Value = "Quant"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "100.42"
$ws.Cells.Item(42, 5).Value = "  -4.24%  "
$ws.Cells.Item(43, 2).Value = "FraxShare"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "5.464"
$ws.Cells.Item(43, 5).Value = "  -2.68%  "
$ws.Cells.Item(44, 2).Value = "TheSandbox"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.3878"
$ws.Cells.Item(44, 5).Value = "  -0.75%  "
$ws.Cells.Item(45, 2).Value = "Aptos"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "6.665"
$ws.Cells.Item(45, 5).Value = "  -3.50%  "
$ws.Cells.Item(46, 2).Value = "Algorand"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.1160"
$ws.Cells.Item(46, 5).Value = "  -2.05%  "
$ws.Cells.Item(47, 2).Value = "Cronos"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.05370"
$ws.Cells.Item(47, 5).Value = "  +0.75%  "
$ws.Cells.Item(48, 2).Value = "Elrond"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "30.43"
$ws.Cells.Item(48, 5).Value = "  -1.57%  "
$ws.Cells.Item(49, 2).Value = "EnergySwap"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "7.668"
$ws.Cells.Item(49, 5).Value = "  -2.28%  "
$ws.Cells.Item(50, 2).Value = "NEARProtocol"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.251"
$ws.Cells.Item(50, 5).Value = "  -1.40%  "
$ws.Cells.Item(51, 2).Value = "Aave"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "51.17"
$ws.Cells.Item(51, 5).Value = "  -0.22%  "
